$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F2").Value = 1.65
$ws.Range("G2").Value = 1.66
$ws.Range("I2").Value = 5.7
$ws.Range("M2").Value = 1.05
$ws.Range("N2").Value = 5.1
$ws.Range("R2").Value = 1.54
$ws.Range("V2").Value = 1.21
$ws.Range("AA2").Value = 140
$ws.Range("AD2").Value = 21
$ws.Range("AH2").Value = 18.5
$ws.Range("AM2").Value = 90
$ws.Range("G3").Value = 120
$ws.Range("K3").Value = 32
$ws.Range("N3").Value = 4.5
$ws.Range("P3").Value = 4.5
$ws.Range("R3").Value = 2.74
$ws.Range("W4").Value = 1.71
$ws.Range("Y4").Value = 1000
$ws.Range("Q5").Value = 1.46
$ws.Range("R5").Value = 1.72
$ws.Range("S5").Value = 2.12
$ws.Range("U5").Value = 2.8
$ws.Range("AB5").Value = 990
$ws.Range("K6").Value = 3.1
$ws.Range("L6").Value = 1.66
$ws.Range("F7").Value = 1.41
$ws.Range("G7").Value = 1.47
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.4
$ws.Range("N7").Value = 4.5
$ws.Range("T7").Value = 1.89
$ws.Range("W7").Value = 3.1
$ws.Range("G8").Value = 2.9
$ws.Range("H8").Value = 2.7
$ws.Range("I8").Value = 2.9
$ws.Range("R8").Value = 1.32
$ws.Range("W8").Value = 1.52
$ws.Range("X8").Value = 16.5
$ws.Range("G9").Value = 1.67
$ws.Range("I9").Value = 11.5
$ws.Range("J9").Value = 3.35
$ws.Range("K9").Value = 3.85
$ws.Range("N9").Value = 2.26
$ws.Range("P9").Value = 1.41
$ws.Range("Q9").Value = 2.96
$ws.Range("R9").Value = 1.14
$ws.Range("S9").Value = 6.4
$ws.Range("T9").Value = 2.88
$ws.Range("U9").Value = 1.42
$ws.Range("W9").Value = 2.52
$ws.Range("AB9").Value = 990
$ws.Range("N10").Value = 2.96
$ws.Range("X10").Value = 9.800000000000001
$ws.Range("AH10").Value = 29
$ws.Range("H11").Value = 36
$ws.Range("I11").Value = 40
$ws.Range("K11").Value = 12
$ws.Range("N11").Value = 6.8
$ws.Range("O11").Value = 1.15
$ws.Range("P11").Value = 2.86
$ws.Range("Q11").Value = 1.48
$ws.Range("R11").Value = 1.71
$ws.Range("S11").Value = 2.24
$ws.Range("T11").Value = 2.8
$ws.Range("U11").Value = 1.54
$ws.Range("X11").Value = 48
$ws.Range("Z11").Value = 420
$ws.Range("AB11").Value = 12
$ws.Range("AC11").Value = 34
$ws.Range("AG11").Value = 15
$ws.Range("AK11").Value = 15
$ws.Range("AL11").Value = 65
$ws.Range("AM11").Value = 520
$ws.Range("AN11").Value = 3.3
